$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterFields")

# The WRITE1SET_REG.BIT0 field (row 13) was removed from the doc -
# clear that row entirely.
$ws.Range("A13:M13").ClearContents()

# The LOCK_TEST_REG register (row 14) no longer documents a fixed
# type / reset value / description - clear those columns but keep
# the register name (A14) and address (C14).
$ws.Range("H14:J14").ClearContents()

# Update the current selection left in the sheet view.
$ws.Range("M20").Select()

$wb.Save()
